# Add one history of CNN reference (row 9 / sheet row 10)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference number
$ws.Range("A10").Value = 9

# Write the URL before the citation text so the shared-string table
# picks up the same ordering as the authored workbook (URL then text).
$ws.Range("D10").Value = "https://blog.athelas.com/a-brief-history-of-cnns-in-image-segmentation-from-r-cnn-to-mask-r-cnn-34ea83205de4"
$ws.Range("C10").Value = 'Dhruv Parthasarathy, "A Brief History of CNNs in Image Segmentation: From R-CNN to Mask R-CNN"'

# Turn the URL into a real hyperlink (adds the external relationship).
$ws.Hyperlinks.Add($ws.Range("D10"), "https://blog.athelas.com/a-brief-history-of-cnns-in-image-segmentation-from-r-cnn-to-mask-r-cnn-34ea83205de4")

# Move the active selection the way the author's session left it.
$ws.Range("C11").Select()
